$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) retains its Text format so values such as
# "69.153.30" or "0.0000284" are stored as literal strings, matching
# the source data, rather than being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.153.30'
$ws.Range('E2').Value = '  -2.35%  '
$ws.Range('D3').Value = '3.460.37'
$ws.Range('E3').Value = '  -4.69%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '576.84'
$ws.Range('E5').Value = '  -4.67%  '
$ws.Range('D6').Value = '191.20'
$ws.Range('E6').Value = '  -4.26%  '
$ws.Range('D7').Value = '0.608'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('D8').Value = '3.450.73'
$ws.Range('E8').Value = '  -4.68%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -6.76%  '
$ws.Range('D11').Value = '0.617'
$ws.Range('E11').Value = '  -4.47%  '
$ws.Range('D12').Value = '51.31'
$ws.Range('E12').Value = '  -4.57%  '
$ws.Range('D13').Value = '0.0000284'
$ws.Range('E13').Value = '  -7.24%  '
$ws.Range('D14').Value = '9.09'
$ws.Range('E14').Value = '  -4.91%  '
$ws.Range('D15').Value = '4.003.76'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').Value = '638.33'
$ws.Range('E16').Value = '  +3.33%  '
$ws.Range('D17').Value = '68.992.58'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '3.446.46'
$ws.Range('E18').Value = '  -5.29%  '
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('D20').Value = '12.23'
$ws.Range('E20').Value = '  -6.03%  '
$ws.Range('D21').Value = '18.06'
$ws.Range('E21').Value = '  -5.20%  '
$ws.Range('D22').Value = '0.939'
$ws.Range('E22').Value = '  -6.03%  '
$ws.Range('D23').Value = '17.80'
$ws.Range('E23').Value = '  -2.66%  '
$ws.Range('D24').Value = '5.28'
$ws.Range('E24').Value = '  -0.92%  '
$ws.Range('D25').Value = '98.48'
$ws.Range('E25').Value = '  -5.06%  '
$ws.Range('D26').Value = '4.26'
$ws.Range('E26').Value = '  -8.03%  '
$ws.Range('D27').Value = '2.85'
$ws.Range('E27').Value = '  -5.22%  '
$ws.Range('D28').Value = '9.85'
$ws.Range('E28').Value = '  -7.08%  '
$ws.Range('D29').Value = '9.22'
$ws.Range('E29').Value = '  -5.38%  '
$ws.Range('D30').Value = '32.21'
$ws.Range('E30').Value = '  -4.44%  '
$ws.Range('D31').Value = '4.31'
$ws.Range('E31').Value = '  -9.03%  '
$ws.Range('D32').Value = '6.68'
$ws.Range('E32').Value = '  -7.22%  '
$ws.Range('D33').Value = '11.54'
$ws.Range('E33').Value = '  -5.77%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  -6.57%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '60.96'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '3.665.16'
$ws.Range('E37').Value = '  -8.12%  '
$ws.Range('D38').Value = '0.0₃0790'
$ws.Range('E38').Value = '  -10.22%  '
$ws.Range('D39').Value = '500.62'
$ws.Range('E39').Value = '  -3.02%  '
$ws.Range('D40').Value = '2.93'
$ws.Range('E40').Value = '  -3.61%  '
$ws.Range('D41').Value = '3.47'
$ws.Range('E41').Value = '  -2.32%  '
$ws.Range('D42').Value = '0.369'
$ws.Range('E42').Value = '  -5.50%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').Value = '3.53'
$ws.Range('E43').Value = '  +70.54%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '0.132'
$ws.Range('E44').Value = '  -2.89%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '34.22'
$ws.Range('E45').Value = '  -6.77%  '
$ws.Range('D46').Value = '0.0438'
$ws.Range('E46').Value = '  -4.93%  '
$ws.Range('E47').Value = '  -3.90%  '
$ws.Range('D48').Value = '2.80'
$ws.Range('E48').Value = '  -3.99%  '
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').Value = '8.07'
$ws.Range('E51').Value = '  -6.19%  '
